$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10:I10").Value = "ОК"
